$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-26 (A: value, B: timestamp serial)
$data = @(
    @(5529, 45821),
    @(5537, 45821.01041666666),
    @(5462, 45821.02083333334),
    @(5407, 45821.03125),
    @(5333, 45821.04166666666),
    @(5337, 45821.05208333334),
    @(5413, 45821.0625),
    @(5338, 45821.07291666666),
    @(5228, 45821.08333333334),
    @(5319, 45821.09375),
    @(5263, 45821.10416666666),
    @(5264, 45821.11458333334),
    @(5232, 45821.125),
    @(5339, 45821.13541666666),
    @(5317, 45821.14583333334),
    @(5289, 45821.15625),
    @(5291, 45821.16666666666),
    @(5317, 45821.17708333334),
    @(5370, 45821.1875),
    @(5343, 45821.19791666666),
    @(5518, 45821.20833333334),
    @(5470, 45821.21875),
    @(5603, 45821.22916666666),
    @(5670, 45821.23958333334),
    @(5802, 45821.25)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

# Delete rows 27 through 41 (old trailing rows no longer present)
$ws.Rows("27:41").Delete() | Out-Null
